$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ A='Sport & Enterteinment'; B=8; C='Movies'; D='adversarial learning'; E='As a movie studio executive, I want to employ adversarial learning techniques to detect deepfake videos that maliciously use our actors'' likenesses, so that we can protect our brand and actors'' reputations.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Music'; D='adversarial learning'; E='As a music data scientist, I want to apply adversarial learning techniques to enhance the robustness of my music genre classification model, ensuring it accurately identifies genre boundaries even in the presence of adversarial attacks.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Movies'; D='cnn'; E='As a visual effects supervisor, I aim to use CNNs to automate the process of identifying and correcting green screen errors in movie footage, reducing production time and costs.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Music'; D='cnn'; E='As a music audio engineer, I aim to utilize CNNs to develop a real-time music event detection system that can identify and timestamp specific musical events (e.g., beats, notes, transitions) in audio recordings.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Movies'; D='conversational agent'; E='As a casual moviegoer, I want a conversational agent to provide summaries and reviews for new movie releases, so that I can decide which films to watch next.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Music'; D='conversational agent'; E='As a musician, I need a conversational AI assistant that can analyze my compositions using machine learning techniques, providing constructive feedback on musical structure, harmony, and creativity.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Movies'; D='decision tree'; E='As a movie studio executive, I want to employ decision trees to predict the financial success of a movie based on its genre, director, and release date, so that I can make informed greenlighting decisions.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Music'; D='decision tree'; E='As a music event organizer, I need a decision tree-based system that recommends suitable artists and bands for events based on their musical style, popularity, and audience preferences extracted from past event data.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Movies'; D='document classification'; E='s a streaming platform administrator, I want to employ document classification to classify user-generated movie summaries into predefined genres (e.g., action, comedy, drama) automatically, so that I can enhance movie recommendation algorithms.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Music'; D='document classification'; E='As a music journalist, I need a document classification system to organize news articles and reviews into topics like album releases, artist interviews, and concert reviews, making it easier to access and reference past work.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Movies'; D='entity extraction'; E='As a movie metadata curator, I want to develop an entity extraction system to automatically identify and extract key entities such as movie titles, directors, actors, and release dates from unstructured movie reviews and articles, so that I can maintain an up-to-date and accurate movie database.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Music'; D='entity extraction'; E='As a music copyright manager, I aim to employ entity extraction to detect and flag copyrighted content in user-uploaded lyrics and compositions, ensuring compliance with intellectual property laws.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Movies'; D='feature selection'; E='As a movie critic, I want to perform feature selection to identify the most influential factors (such as genre, director reputation, and cast popularity) that contribute to a movie''s box office success, so that I can provide deeper insights into what makes a movie commercially appealing.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Music'; D='feature selection'; E='As a music content curator, I want to use feature selection to automate the process of identifying distinctive musical characteristics (e.g., instrumentation, vocal style) that define niche subgenres, facilitating targeted content curation for diverse listener preferences.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Movies'; D='imbalanced dataset'; E='As a researcher studying movie box office performance, I want to handle class imbalance in datasets categorizing movies into blockbuster and non-blockbuster categories, so that I can develop predictive models that accurately forecast commercial success.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Music'; D='imbalanced dataset'; E='As a music recommendation system engineer, I aim to mitigate the impact of imbalanced data by employing techniques that prioritize user feedback on less popular songs, enhancing the diversity and fairness of music recommendations.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Movies'; D='keyword extraction'; E='As a movie critic, I want to develop a keyword extraction system to automatically identify and extract key themes and topics from movie reviews, so that I can categorize and analyze critical opinions effectively.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Music'; D='keyword extraction'; E='As a music content curator, I need a keyword extraction tool that utilizes machine learning to automatically identify prominent themes and topics from music reviews and interviews, aiding in the creation of focused artist profiles and music playlists.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Movies'; D='k-nearest neighbor'; E='As a movie critic, I want to leverage KNN to identify movies with similar stylistic traits and narrative structures, so that I can analyze trends in filmmaking techniques and storytelling approaches over time.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Music'; D='k-nearest neighbor'; E='As a music event organizer, I want to use a k-NN model to recommend suitable opening acts for headline artists based on their musical styles and fan demographics, ensuring cohesive and engaging concert lineups.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Movies'; D='multi-label classification'; E='As a movie content moderator, I want to implement multi-label classification to automatically flag movies with appropriate content warnings (such as violence, nudity, and language) based on scene analysis and dialogue, so that I can ensure viewer suitability.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Music'; D='multi-label classification'; E='As a music recommendation system developer, I aim to build a multi-label classification model that can accurately tag songs with multiple genres (e.g., rock, pop, indie) based on their audio features and lyrical content, enhancing the diversity and relevance of music recommendations.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Movies'; D='neural network'; E='As a special effects supervisor, I want to develop a neural network model to simulate realistic visual effects (such as explosions and CGI creatures) in movies based on physical parameters and environmental conditions, so that I can achieve cinematic realism.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Music'; D='neural network'; E='As a music audio quality evaluator, I want to develop a neural network model that assesses the audio quality of music recordings by analyzing various audio features such as clarity, dynamic range, and noise levels, ensuring high fidelity in music production and distribution.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Movies'; D='random forest'; E='As a movie critic, I want to build a random forest model to predict movie box office success based on features such as genre, cast popularity, and promotional budget, so that I can analyze factors contributing to financial performance.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Music'; D='random forest'; E='As a music content curator, I seek to utilize a random forest model to automatically classify newly released songs into thematic categories (e.g., party music, workout tunes, relaxing melodies) based on their audio attributes, facilitating targeted playlist curation and content discovery.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Movies'; D='semantic similarity'; E='As a movie review aggregator, I want to develop a semantic similarity model to identify and group similar movie reviews based on their thematic content and sentiment, so that I can summarize overall audience reactions more effectively.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Music'; D='semantic similarity'; E='As a music content aggregator, I seek to employ semantic similarity algorithms to cluster music news articles and reviews based on shared topics, enabling efficient content aggregation and providing comprehensive coverage of trending music topics.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Movies'; D='sentiment analysis'; E='As a movie critic, I want to develop a sentiment analysis model to automatically classify movie reviews as positive, neutral, or negative based on textual content, so that I can gauge overall audience sentiment towards a movie.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Music'; D='sentiment analysis'; E='As a music marketing strategist, I want to use sentiment analysis to monitor social media sentiment around new music releases and artist announcements, gauging public reception and adjusting marketing campaigns accordingly.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Movies'; D='speech to text'; E='As a movie review aggregator, I want to implement speech-to-text technology to transcribe audiovisual movie reviews from podcasts and video reviews into text format, so that I can include diverse content sources in my review summaries.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Music'; D='speech to text'; E='As a music journalist, I seek to use speech-to-text tools to transcribe artist interviews and press conferences efficiently, enabling quicker turnaround for articles and features covering music industry news and events.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Movies'; D='text categorization'; E='s a movie script analyst, I want to employ text categorization to classify movie scripts into different narrative styles and thematic categories (e.g., romance, action-adventure, science fiction), so that I can analyze storytelling trends and genre preferences over time.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Music'; D='text categorization'; E='As a music content curator, I need a text categorization model to automatically categorize music reviews and critiques into genres (e.g., rock, jazz, classical) based on textual content, facilitating efficient content tagging and organization.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Movies'; D='unsupervised clustering'; E='As a movie critic, I want to utilize unsupervised clustering to group movies into clusters based on thematic elements and directorial styles, so that I can explore connections between movies and provide insightful reviews on thematic similarities.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Music'; D='unsupervised clustering'; E='As a music playlist curator, I seek to use unsupervised clustering techniques to cluster user-generated playlists into groups based on similarity in musical styles, facilitating the discovery of new playlist ideas and enhancing playlist diversity.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Movies'; D='voice recognition'; E='As a movie subtitle creator, I want to develop a voice recognition system to accurately transcribe movie dialogues from audio tracks into text, so that I can create synchronized subtitles for hearing-impaired viewers.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Music'; D='voice recognition'; E='As a music event organizer, I aim to utilize voice recognition technology for real-time transcription of panel discussions, workshops, and keynote speeches at music conferences and industry events, facilitating knowledge sharing and accessibility for attendees.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Movies'; D='word embedding'; E='As a movie content curator, I want to use word embedding techniques to analyze similarities and relationships between movie plot summaries based on semantic meanings of key words and phrases, so that I can categorize and recommend related movies effectively.'; F='ReAdjusted_CoTPrompt' },
    @{ A='Sport & Enterteinment'; B=8; C='Music'; D='word embedding'; E='As a music content curator, I seek to employ word embedding algorithms to analyze and categorize music-related textual content (e.g., artist bios, album reviews) based on underlying semantic similarities, facilitating more intuitive content organization and retrieval.'; F='ReAdjusted_CoTPrompt' },
)

$startRow = 582
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 3)).Interior.Color = 13223074
    $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 3)).Interior.PatternColor = 13223074
}

$ws.Range("F581:F621").Select()
$excel.ActiveWindow.ScrollRow = 598